$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tidy")

# Update the shared text "Nacionales de Abulón" -> "Productores Nacionales de Abulón"
# This string is used in B16:B20, so update every occurrence to the same new
# text so the workbook keeps a single shared-string entry (edited in place).
$ws.Range("B16:B20").Value2 = "Productores Nacionales de Abulón"

# Update the active cell selection on the "tidy" sheet from F16 to C17
$ws.Range("C17").Select()
